# Added implementation of MSM measure.
# This updates two sheets that were auto-generated from source code analysis:
#   - interfaceOperations (sheet2): operations for com.macro.mall.handler.SwaggerHandler,
#     now also including inherited java.lang.Object operations.
#   - methodNumberOfLines (sheet11): per-method line counts, now also including the
#     previously-missing no-arg constructors for IgnoreUrlsConfig and
#     IgnoreUrlsRemoveJwtFilter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# interfaceOperations
# ---------------------------------------------------------------------------
$wsOps = $wb.Worksheets.Item("interfaceOperations")

$opsRows = @(
    @("com.macro.mall.handler.SwaggerHandler", "equals(java.lang.Object)", "public", "boolean"),
    @("com.macro.mall.handler.SwaggerHandler", "toString()", "public", "java.lang.String"),
    @("com.macro.mall.handler.SwaggerHandler", "swaggerResources()", "public", "reactor.core.publisher.Mono"),
    @("com.macro.mall.handler.SwaggerHandler", "getClass()", "public", "java.lang.Class"),
    @("com.macro.mall.handler.SwaggerHandler", "notifyAll()", "public", "void"),
    @("com.macro.mall.handler.SwaggerHandler", "hashCode()", "public", "int"),
    @("com.macro.mall.handler.SwaggerHandler", "wait()", "public", "void"),
    @("com.macro.mall.handler.SwaggerHandler", "uiConfiguration()", "public", "reactor.core.publisher.Mono"),
    @("com.macro.mall.handler.SwaggerHandler", "securityConfiguration()", "public", "reactor.core.publisher.Mono"),
    @("com.macro.mall.handler.SwaggerHandler", "notify()", "public", "void"),
    @("com.macro.mall.handler.SwaggerHandler", "wait(long)", "public", "void"),
    @("com.macro.mall.handler.SwaggerHandler", "wait(long, int)", "public", "void"),
    @("com.macro.mall.handler.SwaggerHandler", "SwaggerHandler(springfox.documentation.swagger.web.SwaggerResourcesProvider)", "public", "void")
)

for ($i = 0; $i -lt $opsRows.Length; $i++) {
    $r = $i + 2
    $row = $opsRows[$i]
    $wsOps.Range("A$r").Value = $row[0]
    $wsOps.Range("B$r").Value = $row[1]
    $wsOps.Range("C$r").Value = $row[2]
    $wsOps.Range("D$r").Value = $row[3]
}

# ---------------------------------------------------------------------------
# methodNumberOfLines
# ---------------------------------------------------------------------------
$wsLines = $wb.Worksheets.Item("methodNumberOfLines")

$linesRows = @(
    @("com.macro.mall.config.SwaggerResourceConfig", "get()", "7"),
    @("com.macro.mall.config.SwaggerResourceConfig", "swaggerResource(java.lang.String, java.lang.String)", "8"),
    @("com.macro.mall.config.SwaggerResourceConfig", "SwaggerResourceConfig(org.springframework.cloud.gateway.route.RouteLocator, org.springframework.cloud.gateway.config.GatewayProperties)", "4"),
    @("com.macro.mall.config.SwaggerResourceConfig", "get()", "7"),
    @("com.macro.mall.config.SwaggerResourceConfig", "lambda`$get`$4(java.util.List, org.springframework.cloud.gateway.route.RouteDefinition)", "1"),
    @("com.macro.mall.config.SwaggerResourceConfig", "lambda`$get`$3(java.util.List, org.springframework.cloud.gateway.route.RouteDefinition, org.springframework.cloud.gateway.handler.predicate.PredicateDefinition)", "1"),
    @("com.macro.mall.config.SwaggerResourceConfig", "lambda`$get`$2(org.springframework.cloud.gateway.handler.predicate.PredicateDefinition)", "1"),
    @("com.macro.mall.config.SwaggerResourceConfig", "lambda`$get`$1(java.util.List, org.springframework.cloud.gateway.route.RouteDefinition)", "1"),
    @("com.macro.mall.config.SwaggerResourceConfig", "lambda`$get`$0(java.util.List, org.springframework.cloud.gateway.route.Route)", "1"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "IgnoreUrlsConfig()", "1"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "getUrls()", "3"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "setUrls(java.util.List)", "3"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "toString()", "3"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "equals(java.lang.Object)", "23"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "canEqual(java.lang.Object)", "3"),
    @("com.macro.mall.config.IgnoreUrlsConfig", "hashCode()", "7"),
    @("com.macro.mall.handler.SwaggerHandler", "SwaggerHandler(springfox.documentation.swagger.web.SwaggerResourcesProvider)", "3"),
    @("com.macro.mall.handler.SwaggerHandler", "securityConfiguration()", "3"),
    @("com.macro.mall.handler.SwaggerHandler", "uiConfiguration()", "3"),
    @("com.macro.mall.handler.SwaggerHandler", "swaggerResources()", "3"),
    @("com.macro.mall.filter.AuthGlobalFilter", "AuthGlobalFilter()", "1"),
    @("com.macro.mall.filter.AuthGlobalFilter", "filter(org.springframework.web.server.ServerWebExchange, org.springframework.cloud.gateway.filter.GatewayFilterChain)", "18"),
    @("com.macro.mall.filter.AuthGlobalFilter", "getOrder()", "3"),
    @("com.macro.mall.config.ResourceServerConfig", "springSecurityFilterChain(org.springframework.security.config.web.server.ServerHttpSecurity)", "7"),
    @("com.macro.mall.config.ResourceServerConfig", "jwtAuthenticationConverter()", "8"),
    @("com.macro.mall.config.ResourceServerConfig", "ResourceServerConfig(com.macro.mall.authorization.AuthorizationManager, com.macro.mall.config.IgnoreUrlsConfig, com.macro.mall.component.RestfulAccessDeniedHandler, com.macro.mall.component.RestAuthenticationEntryPoint, com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter)", "7"),
    @("com.macro.mall.config.GlobalCorsConfig", "corsFilter()", "10"),
    @("com.macro.mall.MallGatewayApplication", "MallGatewayApplication()", "1"),
    @("com.macro.mall.MallGatewayApplication", "main(java.lang.String[])", "3"),
    @("com.macro.mall.authorization.AuthorizationManager", "AuthorizationManager()", "1"),
    @("com.macro.mall.authorization.AuthorizationManager", "check(reactor.core.publisher.Mono, org.springframework.security.web.server.authorization.AuthorizationContext)", "54"),
    @("com.macro.mall.authorization.AuthorizationManager", "check(reactor.core.publisher.Mono, java.lang.Object)", "1"),
    @("com.macro.mall.authorization.AuthorizationManager", "lambda`$check`$0(java.lang.String)", "1"),
    @("com.macro.mall.component.RestAuthenticationEntryPoint", "RestAuthenticationEntryPoint()", "1"),
    @("com.macro.mall.component.RestAuthenticationEntryPoint", "commence(org.springframework.web.server.ServerWebExchange, org.springframework.security.core.AuthenticationException)", "10"),
    @("com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter", "IgnoreUrlsRemoveJwtFilter()", "1"),
    @("com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter", "filter(org.springframework.web.server.ServerWebExchange, org.springframework.web.server.WebFilterChain)", "14"),
    @("com.macro.mall.component.RestfulAccessDeniedHandler", "RestfulAccessDeniedHandler()", "1"),
    @("com.macro.mall.component.RestfulAccessDeniedHandler", "handle(org.springframework.web.server.ServerWebExchange, org.springframework.security.access.AccessDeniedException)", "10")
)

for ($i = 0; $i -lt $linesRows.Length; $i++) {
    $r = $i + 2
    $row = $linesRows[$i]
    $wsLines.Range("A$r").Value = $row[0]
    $wsLines.Range("B$r").Value = $row[1]
    # "Number of Lines" is stored as text (shared string) in the source data,
    # not as a number -- force text formatting so it round-trips as a string.
    $wsLines.Range("C$r").NumberFormat = "@"
    $wsLines.Range("C$r").Value = $row[2]
}
